$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 860.4545000000001
$ws.Range("I28").Value = 845.9
$ws.Range("J28").Value = 1006
$ws.Range("K28").Value = 845.9
$ws.Range("L28").Value = 1006
$ws.Range("M28").Value = -360.9
$ws.Range("N28").Value = -1976
$ws.Range("H32").Value = 37499
$ws.Range("J32").Value = 39998.332
$ws.Range("L32").Value = 39998.332
$ws.Range("N32").Value = -40650.332
$ws.Range("H43").Value = 5718.375
$ws.Range("I43").Value = 2650
$ws.Range("J43").Value = 7559.4
$ws.Range("K43").Value = 2650
$ws.Range("L43").Value = 7559.4
$ws.Range("M43").Value = -2581
$ws.Range("N43").Value = -7697.4
$ws.Range("H100").Value = 4065.4
$ws.Range("I100").Value = 4206
$ws.Range("J100").Value = 3503
$ws.Range("K100").Value = 4206
$ws.Range("L100").Value = 3503
$ws.Range("M100").Value = -3665
$ws.Range("N100").Value = -4585
$ws.Range("H125").Value = 2749.5715
$ws.Range("I125").Value = 2869.875
$ws.Range("J125").Value = 2589.1667
$ws.Range("K125").Value = 25828.875
$ws.Range("L125").Value = 23302.5003
$ws.Range("M125").Value = -23368.875
$ws.Range("N125").Value = -28222.5003
$ws.Range("H132").Value = 2740.697
$ws.Range("I132").Value = 2394.6206
$ws.Range("K132").Value = 7183.861800000001
$ws.Range("M132").Value = -4653.861800000001
$ws.Range("H137").Value = 3154.5625
$ws.Range("J137").Value = 5249.8
$ws.Range("L137").Value = 15749.4
$ws.Range("N137").Value = -20849.4
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12200915
$ws.Range("I32").Value = 12200915
$ws.Range("K32").Value = 12200915
$ws.Range("M32").Value = -12200628
$ws.Range("H61").Value = 21790330
$ws.Range("I61").Value = 62503388
$ws.Range("K61").Value = 62503388
$ws.Range("M61").Value = -62503176
$ws.Range("H74").Value = 20849570
$ws.Range("I74").Value = 125000950
$ws.Range("J74").Value = 19294.7
$ws.Range("K74").Value = 125000950
$ws.Range("L74").Value = 19294.7
$ws.Range("M74").Value = -125000076
$ws.Range("N74").Value = -21042.7
$ws.Range("H77").Value = 20849570
$ws.Range("I77").Value = 125000950
$ws.Range("J77").Value = 19294.7
$ws.Range("K77").Value = 625004750
$ws.Range("L77").Value = 96473.5
$ws.Range("M77").Value = -625000382
$ws.Range("N77").Value = -105209.5
$ws.Range("H88").Value = 2668.5454
$ws.Range("I88").Value = 2049.75
$ws.Range("J88").Value = 3022.1428
$ws.Range("K88").Value = 2049.75
$ws.Range("L88").Value = 3022.1428
$ws.Range("M88").Value = -1643.75
$ws.Range("N88").Value = -3834.1428
$ws.Range("H91").Value = 2668.5454
$ws.Range("I91").Value = 2049.75
$ws.Range("J91").Value = 3022.1428
$ws.Range("K91").Value = 2049.75
$ws.Range("L91").Value = 3022.1428
$ws.Range("M91").Value = -645.75
$ws.Range("N91").Value = -5830.1428
$ws.Range("H122").Value = 1382.2
$ws.Range("I122").Value = 1137.3334
$ws.Range("K122").Value = 3412.0002
$ws.Range("M122").Value = -962.0001999999999
$ws.Range("H126").Value = 6230
$ws.Range("I126").Value = 6230
$ws.Range("K126").Value = 18690
$ws.Range("M126").Value = -16220
$ws.Range("H132").Value = 6265.4185
$ws.Range("I132").Value = 2854.6667
$ws.Range("J132").Value = 14136.385
$ws.Range("K132").Value = 8564.000100000001
$ws.Range("L132").Value = 42409.155
$ws.Range("M132").Value = -6034.000100000001
$ws.Range("N132").Value = -47469.155
$ws.Range("H136").Value = 21790330
$ws.Range("I136").Value = 62503388
$ws.Range("K136").Value = 187510164
$ws.Range("M136").Value = -187507614
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 355
$ws.Range("J80").Value = 344.25
$ws.Range("L80").Value = 344.25
$ws.Range("N80").Value = -2340.25
$ws.Range("H83").Value = 355
$ws.Range("J83").Value = 344.25
$ws.Range("L83").Value = 1721.25
$ws.Range("N83").Value = -11705.25
$ws.Range("H107").Value = 2131
$ws.Range("J107").Value = 2333.3333
$ws.Range("L107").Value = 2333.3333
$ws.Range("N107").Value = -6173.3333
$ws.Range("H110").Value = 128990
$ws.Range("J110").Value = 128990
$ws.Range("L110").Value = 128990
$ws.Range("N110").Value = -137170
$ws.Range("H134").Value = 36041.332
$ws.Range("I134").Value = 1578.4348
$ws.Range("K134").Value = 4735.3044
$ws.Range("M134").Value = -2200.3044
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 44636.09
$ws.Range("I51").Value = 35124.625
$ws.Range("K51").Value = 35124.625
$ws.Range("M51").Value = -34388.625
$ws.Range("H53").Value = 54558
$ws.Range("J53").Value = 54558
$ws.Range("L53").Value = 54558
$ws.Range("N53").Value = -55772
$ws.Range("H60").Value = 71631.05499999999
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 71631.05499999999
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 71631.05499999999
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -72653.05499999999
$ws.Range("H61").Value = 44636.09
$ws.Range("I61").Value = 35124.625
$ws.Range("K61").Value = 35124.625
$ws.Range("M61").Value = -34776.625
$ws.Range("H68").Value = 50000
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 50000
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H86").Value = 2918.1667
$ws.Range("I86").Value = 2877.5
$ws.Range("K86").Value = 2877.5
$ws.Range("M86").Value = -1754.5
$ws.Range("H89").Value = 2918.1667
$ws.Range("I89").Value = 2877.5
$ws.Range("K89").Value = 2877.5
$ws.Range("M89").Value = -8771.5
$ws.Range("H99").Value = 2574.7144
$ws.Range("I99").Value = 2574.7144
$ws.Range("K99").Value = 2574.7144
$ws.Range("M99").Value = -1076.7144
$ws.Range("H122").Value = 3530.182
$ws.Range("I122").Value = 3483.2
$ws.Range("K122").Value = 10449.6
$ws.Range("M122").Value = -7999.599999999999
$ws.Range("H126").Value = 2574.7144
$ws.Range("I126").Value = 2574.7144
$ws.Range("K126").Value = 7724.1432
$ws.Range("M126").Value = -5254.1432
$ws.Range("H134").Value = 241182.47
$ws.Range("I134").Value = 335029.34
$ws.Range("K134").Value = 1005088.02
$ws.Range("M134").Value = -1002553.02
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4203773.5
$ws.Range("I4").Value = 4480075.5
$ws.Range("K4").Value = 13440226.5
$ws.Range("M4").Value = -13440114.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4950
$ws.Range("J70").Value = 5000
$ws.Range("L70").Value = 5000
$ws.Range("N70").Value = -5540
$ws.Range("H73").Value = 4950
$ws.Range("J73").Value = 5000
$ws.Range("L73").Value = 5000
$ws.Range("N73").Value = -6872
$ws.Range("H93").Value = 59970
$ws.Range("J93").Value = 59970
$ws.Range("L93").Value = 59970
$ws.Range("N93").Value = -63714
$ws.Range("H107").Value = 1538.9
$ws.Range("I107").Value = 1432.1111
$ws.Range("K107").Value = 1432.1111
$ws.Range("M107").Value = 487.8888999999999
$ws.Range("H113").Value = 4140.8335
$ws.Range("I113").Value = 3847.4443
$ws.Range("K113").Value = 3847.4443
$ws.Range("M113").Value = -1677.4443
$ws.Range("H132").Value = 111115460
$ws.Range("I132").Value = 166670690
$ws.Range("K132").Value = 500012070
$ws.Range("M132").Value = -500009540
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3891.3333
$ws.Range("I40").Value = 2385.8572
$ws.Range("K40").Value = 2385.8572
$ws.Range("M40").Value = -2249.8572
$ws.Range("H98").Value = 86441.5
$ws.Range("J98").Value = 86441.5
$ws.Range("L98").Value = 86441.5
$ws.Range("N98").Value = -92431.5
$ws.Range("H100").Value = 5000
$ws.Range("I100").Value = 5000
$ws.Range("K100").Value = 5000
$ws.Range("M100").Value = -4459
$ws.Range("H132").Value = 126774.44
$ws.Range("I132").Value = 78337.38
$ws.Range("K132").Value = 235012.14
$ws.Range("M132").Value = -232482.14
$ws.Range("H136").Value = 48817.04
$ws.Range("I136").Value = 4279.3
$ws.Range("K136").Value = 12837.9
$ws.Range("M136").Value = -10287.9
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H122").Value = 6349.3
$ws.Range("I122").Value = 4370.2
$ws.Range("K122").Value = 13110.6
$ws.Range("M122").Value = -10660.6
$ws.Range("H126").Value = 1352
$ws.Range("I126").Value = 1352
$ws.Range("K126").Value = 4056
$ws.Range("M126").Value = -1586
$ws.Range("H132").Value = 1200.2727
$ws.Range("I132").Value = 1175.2812
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 3525.8436
$ws.Range("L132").Value = 2000
$ws.Range("M132").Value = -995.8435999999997
$ws.Range("N132").Value = -11060
